$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename first column header from "season" to "year"
$ws.Range("A1").Value = "year"

# Insert a new column before the existing "pts_per_game" column (currently column F)
# so the new "pf_per_game" column sits between "team" and "pts_per_game"
$ws.Columns.Item(6).Insert()

# New header for the inserted column, matching the bold header style used by the others
$ws.Range("F1").Value = "pf_per_game"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# New "pf_per_game" values (column F), rows 2..11
$pfValues = @(2.1, 3.3, 2.9, 3.4, 2.2, 3.1, 2.1, 3.2, 2.3, 2.8)
for ($i = 0; $i -lt $pfValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 6).Value = $pfValues[$i]
}

# New "season_y" column appended at the end (column H), header styled like the others
$ws.Range("H1").Value = "season_y"
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Copy the "year" values into the new "season_y" column. Using Copy/PasteSpecial
# (instead of a plain Value assignment) preserves the text nature of the values
# instead of Excel auto-converting the numeric-looking strings into numbers.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Copy()
    $ws.Cells.Item($row, 8).PasteSpecial(-4163)  # xlPasteValues
}
